# Update existing cell values on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 1.57

# Row 4
$ws.Range("G4").Value = 1.86

# Row 5
$ws.Range("N5").Value = 23
$ws.Range("O5").Value = 1.1
$ws.Range("P5").Value = 7
$ws.Range("Q5").Value = 1.33
$ws.Range("R5").Value = 3.4
$ws.Range("S5").Value = 1.11
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 10
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 7
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 23
$ws.Range("AE5").Value = 29
$ws.Range("AF5").Value = 67
$ws.Range("AH5").Value = 41
$ws.Range("AL5").Value = 81
$ws.Range("AM5").Value = 67
$ws.Range("AN5").Value = 3.75

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.38

# Add new row 7 with a full new match record
$row7 = @{
    A  = "vw4ebel5"
    B  = "23/10/2024"
    C  = "20:30"
    D  = "USA - USL CHAMPIONSHIP"
    E  = "Tampa Bay"
    F  = "Hartford Athletic"
    G  = 1.62
    H  = 3.8
    I  = 4.5
    J  = 2.18
    K  = 2.27
    L  = 4.75
    M  = 1.04
    N  = 8.25
    O  = 1.23
    P  = 3.8
    Q  = 1.7
    R  = 2.07
    S  = 1.34
    T  = 3
    U  = 1.72
    V  = 2.02
    W  = 7.9
    X  = 8.25
    Y  = 8
    Z  = 12.5
    AA = 12.5
    AB = 23
    AC = 8.25
    AD = 7.7
    AE = 15
    AF = 60
    AG = 400
    AH = 14.5
    AI = 28
    AJ = 15
    AK = 80
    AL = 40
    AM = 40
    AN = 3.6
    AO = 7.9
    AP = 16
    AQ = 25
    AR = 50
    AS = 200
    AT = 3
    AU = 7.3
    AV = 60
    AW = 6.4
    AX = 25
    AY = 29
    AZ = 150
    BA = 150
    BB = 350
    BC = 51
    BD = 51
}

foreach ($col in $row7.Keys) {
    $ws.Range("$col`7").Value = $row7[$col]
}
